# Generate Report for Handoff
# Update the "In Translation" status cells to "Ready for handoff" and bump
# the associated timestamp cells forward, on all three sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status columns (E, F) + HO Xliff generate date (G)
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-02 09:08:18"

# zh-cn sheet: Status column (C) + Latest Handoff Datetime (H)
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-02 09:08:13"

# de-de sheet: Status column (C) + Latest Handoff Datetime (H)
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-02 09:08:18"

# Widen the columns holding the status text so their widths follow the
# now-longer "Ready for handoff" label (matches Excel's own recalculated
# column widths after the cell content changed, ~17.22 characters wide).
$newStatusColWidth = 16.333333333333336
$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusColWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusColWidth
